# Update the "Förändrad" (Changed) date column (C) for all data rows.
# All cells in C2:C230 previously held the serial date value 45180
# (2023-09-11) and are being bumped to 45181 (2023-09-12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 230 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2() -eq 45180) {
        $cell.Value2 = 45181
    }
}
